$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Columns.Item(1).Insert()
$col = $ws.Columns.Item(1)
$col.VerticalAlignment = -4160
for ($i=4; $i -le 26; $i++) {
    $ws.Range("A$i").Value = $i - 3
}
$col.ColumnWidth = 2.6666666666666665
Write-Host "done"
